# Update "想去人数" (want-to-go count) values in column F across sheets.
$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 634
$ws1.Range("F9").Value  = 6467
$ws1.Range("F13").Value = 5052
$ws1.Range("F14").Value = 112
$ws1.Range("F21").Value = 326
$ws1.Range("F27").Value = 1060
$ws1.Range("F30").Value = 594
$ws1.Range("F33").Value = 110
$ws1.Range("F46").Value = 144

# 演出 (Performance) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F26").Value = 331
$ws2.Range("F27").Value = 417

# 本地生活 (Local Life) sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value  = 1534
$ws3.Range("F11").Value = 882

# 全部类型 (All Types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 634
$ws4.Range("F9").Value  = 1534
$ws4.Range("F12").Value = 6467
$ws4.Range("F18").Value = 326
$ws4.Range("F25").Value = 1060
$ws4.Range("F28").Value = 594
$ws4.Range("F38").Value = 331
$ws4.Range("F39").Value = 417
